# Auto-generated edit applying cryptos price/volume refresh (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.325.55"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.274.51"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.89%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "2.617.68"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "2.277.90"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "42.316.01"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("E34").Value = "  +4.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0807"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +21.03%  "
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.90%  "
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0304"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.61%  "
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("E45").Value = "  +6.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.98%  "
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.30%  "
